$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 184, pushing the existing rows 184..250
# down to 185..251 (matches the dimension growing from A1:R250 to A1:R251).
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new record.
$ws.Range("A184").Value = 7
$ws.Range("B184").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C184").Value = "Ñuble"
$ws.Range("D184").Value = 44795
$ws.Range("E184").Value = 16
$ws.Range("F184").Value = 100112043
$ws.Range("G184").Value = "Pepino ensalada"
$ws.Range("H184").Value = "Sin especificar"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 120
$ws.Range("K184").Value = 22000
$ws.Range("L184").Value = 23000
$ws.Range("M184").Value = 22500
$ws.Range("N184").Value = "$/caja 60 unidades"
$ws.Range("O184").Value = "Región de Arica y Parinacota"
$ws.Range("P184").Value = 375
$ws.Range("Q184").Value = 60
$ws.Range("R184").Value = "Hortaliza"
